# Add segment & page.
#
# The paragraph that carries the "_GoBack" bookmark (the last body
# paragraph of the "segment" section, right before the trailing blank
# paragraph / sectPr) gets a new empty paragraph inserted immediately
# before it and another empty paragraph inserted immediately after it.
# Both new paragraphs only carry the same <w:rPr><w:rFonts
# w:hint="eastAsia"/></w:rPr> paragraph mark formatting already used by
# the existing blank paragraphs in this document - no run/text content.

$d = $word.ActiveDocument

$emptyParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p>'

# Locate the paragraph holding the _GoBack bookmark.
$bookmarkRange = $d.Bookmarks("_GoBack").Range
$targetParagraph = $bookmarkRange.Paragraphs(1)
$targetRange = $targetParagraph.Range

# Insert a new, empty paragraph immediately before the target paragraph.
$beforePoint = $d.Range($targetRange.Start, $targetRange.Start)
[void]$beforePoint.InsertXML($emptyParaXml)

# Re-resolve the target paragraph (the document shifted) via the
# bookmark again, then insert a new, empty paragraph immediately after
# it.
$targetRange = $d.Bookmarks("_GoBack").Range.Paragraphs(1).Range
$afterPoint = $d.Range($targetRange.End, $targetRange.End)
[void]$afterPoint.InsertXML($emptyParaXml)

Write-Host "Paragraphs after edit: $($d.Paragraphs.Count)"
